$d = $word.ActiveDocument

# Replace every occurrence of "July 01, 2022" with "July 02, 2022".
# We only touch the two digits that actually change ("01" -> "02")
# so the surrounding run structure is left untouched.
$searchStart = 0
while ($true) {
    $rng = $d.Range($searchStart, $d.Content.End)
    $found = $rng.Find.Execute("July 01, 2022")
    if (-not $found) { break }

    $dayRng = $d.Range($rng.Start + 5, $rng.Start + 7)
    $dayRng.Text = "02"

    $searchStart = $rng.End
}

# Replace "August 30, 2022" with "August 31, 2022" the same way.
$searchStart = 0
while ($true) {
    $rng = $d.Range($searchStart, $d.Content.End)
    $found = $rng.Find.Execute("August 30, 2022")
    if (-not $found) { break }

    $dayRng = $d.Range($rng.Start + 7, $rng.Start + 9)
    $dayRng.Text = "31"

    $searchStart = $rng.End
}
